$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.470157689643315
$ws.Range("C2").Value = 0.2339133357980074
$ws.Range("D2").Value = 0.1075226812808836
$ws.Range("E2").Value = 0.05439836397515307
$ws.Range("F2").Value = 2.312878378896968
$ws.Range("I2").Value = 1.736421287184186
$ws.Range("L2").Value = 0.2243747179935482
$ws.Range("M2").Value = 0.3159931720355544

$ws.Range("B3").Value = 1.386189995712527
$ws.Range("C3").Value = 0.2037510629317296
$ws.Range("D3").Value = 0.1076720869940289
$ws.Range("E3").Value = 0.05394962065764108
$ws.Range("F3").Value = 2.25933059293196
$ws.Range("I3").Value = 1.71430840627103
$ws.Range("L3").Value = 0.2214141773817389
$ws.Range("M3").Value = 0.3028853873112851

$ws.Range("B4").Value = 1.335603170011154
$ws.Range("C4").Value = 0.1852597149889164
$ws.Range("D4").Value = 0.1077912464957294
$ws.Range("E4").Value = 0.05366683351482848
$ws.Range("F4").Value = 2.227720188075978
$ws.Range("I4").Value = 1.701542869023527
$ws.Range("L4").Value = 0.2197126688249611
$ws.Range("M4").Value = 0.2950430642277766

$ws.Range("B5").Value = 1.315231812445006
$ws.Range("C5").Value = 0.1777308649235749
$ws.Range("D5").Value = 0.1078466504850795
$ws.Range("E5").Value = 0.05354975700495768
$ws.Range("F5").Value = 2.215155243583865
$ws.Range("I5").Value = 1.696543617103487
$ws.Range("L5").Value = 0.2190484816087874
$ws.Range("M5").Value = 0.2918989163573613

$ws.Range("B6").Value = 1.311863833676057
$ws.Range("C6").Value = 0.1764810812508131
$ws.Range("D6").Value = 0.1078562622546357
$ws.Range("E6").Value = 0.05353020510589612
$ws.Range("F6").Value = 2.213087895109055
$ws.Range("I6").Value = 1.695725705140944
$ws.Range("L6").Value = 0.218939955574541
$ws.Range("M6").Value = 0.2913799500532264

$ws.Range("B7").Value = 1.335327450952889
$ws.Range("C7").Value = 0.1851581526715051
$ws.Range("D7").Value = 0.1077919660427007
$ws.Range("E7").Value = 0.05366526204123101
$ws.Range("F7").Value = 2.22754945423722
$ws.Range("I7").Value = 1.701474627992269
$ws.Range("L7").Value = 0.219703593229859
$ws.Range("M7").Value = 0.2950004521198721

$ws.Range("B8").Value = 1.441003940587848
$ws.Range("C8").Value = 0.223507103662314
$ws.Range("D8").Value = 0.1075684741533998
$ws.Range("E8").Value = 0.05424513512333995
$ws.Range("F8").Value = 2.294150467946096
$ws.Range("I8").Value = 1.728627446223825
$ws.Range("L8").Value = 0.2233297603924953
$ws.Range("M8").Value = 0.3114307952182003

$ws.Range("B9").Value = 1.655972928618723
$ws.Range("C9").Value = 0.2989665458969455
$ws.Range("D9").Value = 0.1073500978286717
$ws.Range("E9").Value = 0.05532533564858788
$ws.Range("F9").Value = 2.434931384048838
$ws.Range("I9").Value = 1.788380427314394
$ws.Range("L9").Value = 0.2313660915368416
$ws.Range("M9").Value = 0.3452920438299927

$ws.Range("B10").Value = 1.818709298254589
$ws.Range("C10").Value = 0.3546131732749132
$ws.Range("D10").Value = 0.1073269170639364
$ws.Range("E10").Value = 0.05608520853995547
$ws.Range("F10").Value = 2.544738275459451
$ws.Range("I10").Value = 1.83634290636077
$ws.Range("L10").Value = 0.2378395095866921
$ws.Range("M10").Value = 0.371184524829907

$ws.Range("B11").Value = 1.893803776385084
$ws.Range("C11").Value = 0.3799843339745053
$ws.Range("D11").Value = 0.1073468735176775
$ws.Range("E11").Value = 0.05642377703135537
$ws.Range("F11").Value = 2.596114556266656
$ws.Range("I11").Value = 1.85906536543007
$ws.Range("L11").Value = 0.2409091785256265
$ws.Range("M11").Value = 0.3831873330131685

$ws.Range("B12").Value = 1.922394475568751
$ws.Range("C12").Value = 0.3896007724806054
$ws.Range("D12").Value = 0.1073588733540234
$ws.Range("E12").Value = 0.05655097919180019
$ws.Range("F12").Value = 2.615777120314391
$ws.Range("I12").Value = 1.867801379878301
$ws.Range("L12").Value = 0.2420896178969656
$ws.Range("M12").Value = 0.3877649369682601

$ws.Range("B13").Value = 1.916230090374938
$ws.Range("C13").Value = 0.3875292919018989
$ws.Range("D13").Value = 0.1073560904747168
$ws.Range("E13").Value = 0.05652362845149117
$ws.Range("F13").Value = 2.611533165379427
$ws.Range("I13").Value = 1.86591404969495
$ws.Range("L13").Value = 0.2418345865510645
$ws.Range("M13").Value = 0.3867776243353234

$ws.Range("B14").Value = 1.8961528575839
$ws.Range("C14").Value = 0.3807753011143973
$ws.Range("D14").Value = 0.1073477713861308
$ws.Range("E14").Value = 0.05643426211327984
$ws.Range("F14").Value = 2.597728033850444
$ws.Range("I14").Value = 1.859781438224303
$ws.Range("L14").Value = 0.2410059324083846
$ws.Range("M14").Value = 0.3835632849214292

$ws.Range("B15").Value = 1.883875062757284
$ws.Range("C15").Value = 0.3766394707894847
$ws.Range("D15").Value = 0.1073432559706262
$ws.Range("E15").Value = 0.0563793920338167
$ws.Range("F15").Value = 2.589299090884168
$ws.Range("I15").Value = 1.856042208314136
$ws.Range("L15").Value = 0.2405007068275751
$ws.Range("M15").Value = 0.3815986337472665

$ws.Range("B16").Value = 1.813823203404638
$ws.Range("C16").Value = 0.3529563195009473
$ws.Range("D16").Value = 0.1073262314671766
$ws.Range("E16").Value = 0.0560629406532378
$ws.Range("F16").Value = 2.541409604931886
$ws.Range("I16").Value = 1.83487626119026
$ws.Range("L16").Value = 0.2376414178849018
$ws.Range("M16").Value = 0.370404638094314

$ws.Range("B17").Value = 1.771122040365412
$ws.Range("C17").Value = 0.3384426260527107
$ws.Range("D17").Value = 0.1073236398620807
$ws.Range("E17").Value = 0.05586699958285379
$ws.Range("F17").Value = 2.512397686734545
$ws.Range("I17").Value = 1.822124263937056
$ws.Range("L17").Value = 0.235919366698667
$ws.Range("M17").Value = 0.3635950180075511

$ws.Range("B18").Value = 1.746661635153373
$ws.Range("C18").Value = 0.3301000775857119
$ws.Range("D18").Value = 0.1073250167739559
$ws.Range("E18").Value = 0.05575363013559276
$ws.Range("F18").Value = 2.495844727708345
$ws.Range("I18").Value = 1.814874658575263
$ws.Range("L18").Value = 0.2349406385553863
$ws.Range("M18").Value = 0.3596993961436112

$ws.Range("B19").Value = 1.7383969399732
$ws.Range("C19").Value = 0.327276330740176
$ws.Range("D19").Value = 0.1073259737542358
$ws.Range("E19").Value = 0.0557151297841294
$ws.Range("F19").Value = 2.490263102934847
$ws.Range("I19").Value = 1.812434620968745
$ws.Range("L19").Value = 0.2346112741766149
$ws.Range("M19").Value = 0.3583840231585427

$ws.Range("B20").Value = 1.775657277879986
$ws.Range("C20").Value = 0.3399870753639789
$ws.Range("D20").Value = 0.1073236185547088
$ws.Range("E20").Value = 0.05588792701036649
$ws.Range("F20").Value = 2.515472174854693
$ws.Range("I20").Value = 1.823472927132258
$ws.Range("L20").Value = 0.2361014654277653
$ws.Range("M20").Value = 0.3643177296442559

$ws.Range("B21").Value = 1.902045839093944
$ws.Range("C21").Value = 0.3827588639051669
$ws.Range("D21").Value = 0.1073500938750556
$ws.Range("E21").Value = 0.05646053836678711
$ws.Range("F21").Value = 2.601777283543157
$ws.Range("I21").Value = 1.86157915359847
$ws.Range("L21").Value = 0.2412488385667046
$ws.Range("M21").Value = 0.384506533634756

$ws.Range("B22").Value = 1.985546452005565
$ws.Range("C22").Value = 0.4107650891094181
$ws.Range("D22").Value = 0.1073933194186694
$ws.Range("E22").Value = 0.05682891407440493
$ws.Range("F22").Value = 2.659393261895985
$ws.Range("I22").Value = 1.88725103794421
$ws.Range("L22").Value = 0.2447180242865556
$ws.Range("M22").Value = 0.3978900228595279

$ws.Range("B23").Value = 1.940898085229662
$ws.Range("C23").Value = 0.395812611086626
$ws.Range("D23").Value = 0.1073678584008633
$ws.Range("E23").Value = 0.05663283616656489
$ws.Range("F23").Value = 2.628530884024769
$ws.Range("I23").Value = 1.873478764689096
$ws.Range("L23").Value = 0.2428568177432027
$ws.Range("M23").Value = 0.3907296590866949

$ws.Range("B24").Value = 1.773606619659006
$ws.Range("C24").Value = 0.3392888250021429
$ws.Range("D24").Value = 0.1073236192629992
$ws.Range("E24").Value = 0.05587846796509766
$ws.Range("F24").Value = 2.514081805330903
$ws.Range("I24").Value = 1.822862942176357
$ws.Range("L24").Value = 0.2360191033938861
$ws.Range("M24").Value = 0.3639909315497007

$ws.Range("B25").Value = 1.596981479964313
$ws.Range("C25").Value = 0.2785201283110155
$ws.Range("D25").Value = 0.1073853030546488
$ws.Range("E25").Value = 0.05503910944908874
$ws.Range("F25").Value = 2.395739331198484
$ws.Range("I25").Value = 1.771509277352621
$ws.Range("L25").Value = 0.2290924722070571
$ws.Range("M25").Value = 0.2918989163573613
